# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values for each coin row, plus a position swap between the HuobiToken and
# ImmutableX rows (33 <-> 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new D (price) value, written as plain text (these values already
# contain two "." separators, e.g. thousand-grouping like "27.060.05", so
# Excel cannot reinterpret them as a number and they stay text automatically)
$priceUpdatesSafe = @{
    2  = "27.060.05"
    3  = "1.818.33"
    12 = "1.821.53"
    21 = "27.056.05"
    24 = "2.040.98"
}

# row -> new D (price) value that looks like a plain decimal number
# (single "."), so the cell format must be forced to Text first or Excel
# would silently convert it to a numeric value (e.g. "310.60" -> 310.6).
$priceUpdatesText = @{
    4  = "1.003"
    5  = "310.60"
    7  = "0.4225"
    8  = "0.3662"
    9  = "0.07215"
    10 = "0.8445"
    11 = "20.83"
    13 = "6.636"
    14 = "0.07083"
    15 = "5.274"
    16 = "89.30"
    17 = "1.004"
    18 = "0.000008821"
    20 = "14.95"
    22 = "5.106"
    23 = "10.82"
    25 = "1.973"
    26 = "151.75"
    27 = "2.253"
    28 = "18.33"
    29 = "5.200"
    30 = "116.11"
    31 = "0.08786"
    32 = "1.178"
    35 = "4.415"
    36 = "1.001"
    37 = "1.097"
    38 = "0.01964"
    39 = "0.05240"
    40 = "7.288"
    41 = "2.868"
    42 = "0.1690"
    43 = "0.5020"
    44 = "8.567"
    45 = "10.57"
    46 = "0.4741"
    47 = "105.97"
    48 = "1.001"
    49 = "0.06371"
    50 = "1.650"
    51 = "1.871"
}

# row -> new E (volume %) value (without the surrounding spaces/percent sign)
$volUpdates = @{
    2  = "-2.43"
    3  = "-1.68"
    4  = "-1.07"
    5  = "-3.03"
    6  = "-1.01"
    7  = "-2.11"
    8  = "-2.17"
    9  = "-2.11"
    10 = "-3.94"
    11 = "-3.95"
    12 = "-1.57"
    13 = "-1.48"
    14 = "-0.57"
    15 = "-3.39"
    16 = "+1.08"
    17 = "-1.16"
    18 = "-2.04"
    19 = "-1.11"
    20 = "-3.41"
    21 = "-2.46"
    22 = "-2.47"
    23 = "-2.70"
    24 = "-2.13"
    25 = "-2.06"
    26 = "-2.55"
    27 = "+5.08"
    28 = "-1.62"
    29 = "-3.78"
    30 = "-2.47"
    31 = "-1.95"
    32 = "-4.54"
    35 = "-3.44"
    36 = "-1.26"
    37 = "-3.50"
    38 = "-0.50"
    39 = "-2.11"
    40 = "+0.35"
    41 = "-0.46"
    42 = "+0.01"
    43 = "-2.54"
    44 = "-2.92"
    45 = "-1.34"
    46 = "-0.31"
    47 = "-3.44"
    48 = "-1.21"
    49 = "-1.88"
    50 = "-2.86"
    51 = "+0.48"
}

foreach ($row in $priceUpdatesSafe.Keys) {
    $price = $priceUpdatesSafe[$row]
    $ws.Range("D$row").Value = $price
}

foreach ($row in $priceUpdatesText.Keys) {
    $price = $priceUpdatesText[$row]
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $price
}

foreach ($row in $volUpdates.Keys) {
    $pct = $volUpdates[$row]
    $ws.Range("E$row").Value = "  $pct%  "
}

# Rows 33 and 34 swap their coin/link/price/volume data:
# old row33 = ImmutableX, old row34 = HuobiToken
# new row33 = HuobiToken, new row34 = ImmutableX (with refreshed numbers)
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.960"
$ws.Range("E33").Value = "  +1.28%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7381"
$ws.Range("E34").Value = "  -5.28%  "
